$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"

$ws.Range("A2").Value = "shopping with Escher"
$ws.Range("B2").Value = "1951-05-06"
$ws.Range("E2").Value = "Enormous Plastic Bottle"
$ws.Range("F2").Value = "Mustard Seed"
$ws.Range("H2").Value = "3.0"
